$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.57"
$ws.Range("E2").Value = "'-1.98%"
$ws.Range("D3").Value = "'37.83"
$ws.Range("E3").Value = "'-4.24%"
$ws.Range("D4").Value = "'5.059"
$ws.Range("E4").Value = "'-1.36%"
$ws.Range("D5").Value = "'0.07908"
$ws.Range("E5").Value = "'-3.41%"
$ws.Range("D6").Value = "'2.068"
$ws.Range("E6").Value = "'4.55%"
$ws.Range("D7").Value = "'4.398"
$ws.Range("E7").Value = "'4.05%"
$ws.Range("D8").Value = "'8.251"
$ws.Range("E8").Value = "'0.99%"
$ws.Range("E9").Value = "'1.20%"
$ws.Range("D10").Value = "'0.9284"
$ws.Range("E10").Value = "'0.14%"
$ws.Range("D11").Value = "'0.1295"
$ws.Range("E11").Value = "'-8.16%"
$ws.Range("D12").Value = "'0.1910"
$ws.Range("E12").Value = "'-2.43%"
$ws.Range("D13").Value = "'0.08706"
$ws.Range("E13").Value = "'-3.43%"
$ws.Range("D14").Value = "'0.03459"
$ws.Range("E14").Value = "'-1.00%"
$ws.Range("D15").Value = "'0.09759"
$ws.Range("E15").Value = "'-0.56%"
$ws.Range("D16").Value = "'0.001398"
$ws.Range("E16").Value = "'-0.23%"
$ws.Range("D17").Value = "'0.006029"
$ws.Range("E17").Value = "'2.47%"
$ws.Range("D18").Value = "'3.560"
$ws.Range("E18").Value = "'-3.13%"
$ws.Range("D19").Value = "'0.3442"
$ws.Range("E19").Value = "'-0.56%"
$ws.Range("D20").Value = "'0.1302"
$ws.Range("E20").Value = "'-3.64%"
$ws.Range("D21").Value = "'5.023"
$ws.Range("E21").Value = "'8.22%"
$ws.Range("D22").Value = "'0.2522"
$ws.Range("E22").Value = "'4.20%"
$ws.Range("D23").Value = "'0.04344"
$ws.Range("E23").Value = "'-0.67%"
$ws.Range("D24").Value = "'0.001222"
$ws.Range("E24").Value = "'-1.66%"
$ws.Range("D25").Value = "'0.004602"
$ws.Range("E25").Value = "'-4.20%"
$ws.Range("E26").Value = "'176.75%"
$ws.Range("D39").Value = "'0.02257"
$ws.Range("E39").Value = "'4.75%"
$ws.Range("D40").Value = "'0.05072"
$ws.Range("E40").Value = "'-2.64%"
$ws.Range("D41").Value = "'0.007596"
$ws.Range("E41").Value = "'0.51%"
$ws.Range("D42").Value = "'0.009976"
$ws.Range("E42").Value = "'2.29%"
$ws.Range("D43").Value = "'0.1359"
$ws.Range("E43").Value = "'-1.05%"
$ws.Range("D44").Value = "'0.002029"
$ws.Range("E44").Value = "'-3.99%"
$ws.Range("D45").Value = "'0.008842"
$ws.Range("E45").Value = "'-10.29%"
$ws.Range("D46").Value = "'0.00006517"
$ws.Range("E46").Value = "'2.31%"
$ws.Range("E47").Value = "'0.61%"
$ws.Range("D48").Value = "'0.003011"
$ws.Range("E48").Value = "'8.95%"
$ws.Range("D49").Value = "'0.001206"
$ws.Range("E49").Value = "'20.72%"
$ws.Range("D50").Value = "'0.00002110"
$ws.Range("E50").Value = "'0.61%"
$ws.Range("E51").Value = "'0.61%"
